$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the diff. NumberFormat is forced to text ("@")
# before each assignment so that numeric-looking strings (e.g. "325.20",
# "0.000009920", "28.432.95") are preserved exactly as text, matching the
# original inlineStr/text cell contents, rather than being auto-converted
# to numbers by Excel.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.432.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.53%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.27%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.20"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.41%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.09%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4551"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.13%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3825"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.16%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07808"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9846"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.23%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.44"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.64%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.839.21"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.11%  "

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.888"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.72%  "

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.627"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.74%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06915"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.04%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.04%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.43"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.67%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009920"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.01%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.65"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.61%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.11%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.447.17"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.46%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.245"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.58%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.91%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.089"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.71%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.068.68"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.26%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.04"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.27%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.638"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.27%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.22"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.72%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.894"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.88%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09268"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.97%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9023"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.84%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.266"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.02%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.314"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.94%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.290"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.76%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05669"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.85%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.150"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.60%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02035"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.85%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.622"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.04%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5539"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.92%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1762"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.72%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.588"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.71%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07122"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.94%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.48"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.47%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5225"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.80%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.60%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.106"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.27%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.47%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "111.62"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.99%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.431"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.25%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.007"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.08%  "
